$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Update the status text "Ready for handoff" -> "In Translation" everywhere it is used:
# Overview!E2 (zh-cn status), Overview!F2 (de-de status), zh-cn!C2 (Status), de-de!C2 (Status)
$wsOverview.Range("E2").Value = "In Translation"
$wsOverview.Range("F2").Value = "In Translation"
$wsZhCn.Range("C2").Value = "In Translation"
$wsDeDe.Range("C2").Value = "In Translation"

# Narrow the Status columns to fit the shorter text (columns were auto-fit to the
# previous, longer status string and need to shrink to match the new text width).
$wsOverview.Columns.Item(5).ColumnWidth = 12.42  # column E
$wsOverview.Columns.Item(6).ColumnWidth = 12.42  # column F
$wsZhCn.Columns.Item(3).ColumnWidth = 12.42       # column C
$wsDeDe.Columns.Item(3).ColumnWidth = 12.42       # column C
